# Update "想去人数" (F column) values on the "展览", "演出" and "全部类型"
# sheets to reflect the newly generated output numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 620
$ws1.Range("F8").Value = 1163
$ws1.Range("F9").Value = 3907
$ws1.Range("F10").Value = 82

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 52

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 620
$ws4.Range("F8").Value = 1163
$ws4.Range("F9").Value = 3907
$ws4.Range("F10").Value = 82
$ws4.Range("F11").Value = 52
